$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$c = $t.Cell(1, 1)
$rng = $c.Range
$rng.MoveEnd(1, -1)
$rng.Find.Execute("77÷4=19, 1", $true, $false, $false, $false, $false, $true, 1, $false, "55÷3=18, 1", 2) | Out-Null

$c = $t.Cell(1, 2)
$rng = $c.Range
$rng.MoveEnd(1, -1)
$rng.Find.Execute("72÷2=36, 0", $true, $false, $false, $false, $false, $true, 1, $false, "57÷5=11, 2", 2) | Out-Null

$c = $t.Cell(1, 3)
$rng = $c.Range
$rng.MoveEnd(1, -1)
$rng.Find.Execute("71÷6=11, 5", $true, $false, $false, $false, $false, $true, 1, $false, "28÷5=5, 3", 2) | Out-Null

$c = $t.Cell(1, 4)
$rng = $c.Range
$rng.MoveEnd(1, -1)
$rng.Find.Execute("44÷8=5, 4", $true, $false, $false, $false, $false, $true, 1, $false, "23÷8=2, 7", 2) | Out-Null

$c = $t.Cell(1, 5)
$rng = $c.Range
$rng.MoveEnd(1, -1)
$rng.Find.Execute("84÷2=42, 0", $true, $false, $false, $false, $false, $true, 1, $false, "13÷5=2, 3", 2) | Out-Null

$c = $t.Cell(5, 1)
$rng = $c.Range
$rng.MoveEnd(1, -1)
$rng.Find.Execute("40÷9=4, 4", $true, $false, $false, $false, $false, $true, 1, $false, "71÷7=10, 1", 2) | Out-Null

$c = $t.Cell(5, 2)
$rng = $c.Range
$rng.MoveEnd(1, -1)
$rng.Find.Execute("73÷7=10, 3", $true, $false, $false, $false, $false, $true, 1, $false, "65÷2=32, 1", 2) | Out-Null

$c = $t.Cell(5, 3)
$rng = $c.Range
$rng.MoveEnd(1, -1)
$rng.Find.Execute("65÷3=21, 2", $true, $false, $false, $false, $false, $true, 1, $false, "38÷3=12, 2", 2) | Out-Null

$c = $t.Cell(5, 4)
$rng = $c.Range
$rng.MoveEnd(1, -1)
$rng.Find.Execute("21÷9=2, 3", $true, $false, $false, $false, $false, $true, 1, $false, "28÷4=7, 0", 2) | Out-Null

$c = $t.Cell(5, 5)
$rng = $c.Range
$rng.MoveEnd(1, -1)
$rng.Find.Execute("16÷8=2, 0", $true, $false, $false, $false, $false, $true, 1, $false, "80÷2=40, 0", 2) | Out-Null

$c = $t.Cell(9, 1)
$rng = $c.Range
$rng.MoveEnd(1, -1)
$rng.Find.Execute("14÷7=2, 0", $true, $false, $false, $false, $false, $true, 1, $false, "17÷4=4, 1", 2) | Out-Null

$c = $t.Cell(9, 2)
$rng = $c.Range
$rng.MoveEnd(1, -1)
$rng.Find.Execute("33÷6=5, 3", $true, $false, $false, $false, $false, $true, 1, $false, "84÷2=42, 0", 2) | Out-Null

$c = $t.Cell(9, 3)
$rng = $c.Range
$rng.MoveEnd(1, -1)
$rng.Find.Execute("42÷7=6, 0", $true, $false, $false, $false, $false, $true, 1, $false, "19÷2=9, 1", 2) | Out-Null

$c = $t.Cell(9, 4)
$rng = $c.Range
$rng.MoveEnd(1, -1)
$rng.Find.Execute("17÷2=8, 1", $true, $false, $false, $false, $false, $true, 1, $false, "13÷6=2, 1", 2) | Out-Null

$c = $t.Cell(9, 5)
$rng = $c.Range
$rng.MoveEnd(1, -1)
$rng.Find.Execute("21÷9=2, 3", $true, $false, $false, $false, $false, $true, 1, $false, "27÷6=4, 3", 2) | Out-Null

$c = $t.Cell(13, 1)
$rng = $c.Range
$rng.MoveEnd(1, -1)
$rng.Find.Execute("14÷6=2, 2", $true, $false, $false, $false, $false, $true, 1, $false, "98÷7=14, 0", 2) | Out-Null

$c = $t.Cell(13, 2)
$rng = $c.Range
$rng.MoveEnd(1, -1)
$rng.Find.Execute("19÷8=2, 3", $true, $false, $false, $false, $false, $true, 1, $false, "20÷5=4, 0", 2) | Out-Null

$c = $t.Cell(13, 3)
$rng = $c.Range
$rng.MoveEnd(1, -1)
$rng.Find.Execute("92÷8=11, 4", $true, $false, $false, $false, $false, $true, 1, $false, "42÷9=4, 6", 2) | Out-Null

$c = $t.Cell(13, 4)
$rng = $c.Range
$rng.MoveEnd(1, -1)
$rng.Find.Execute("54÷5=10, 4", $true, $false, $false, $false, $false, $true, 1, $false, "19÷6=3, 1", 2) | Out-Null

$c = $t.Cell(13, 5)
$rng = $c.Range
$rng.MoveEnd(1, -1)
$rng.Find.Execute("27÷8=3, 3", $true, $false, $false, $false, $false, $true, 1, $false, "53÷8=6, 5", 2) | Out-Null

$c = $t.Cell(17, 1)
$rng = $c.Range
$rng.MoveEnd(1, -1)
$rng.Find.Execute("28÷7=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "89÷2=44, 1", 2) | Out-Null

$c = $t.Cell(17, 2)
$rng = $c.Range
$rng.MoveEnd(1, -1)
$rng.Find.Execute("34÷7=4, 6", $true, $false, $false, $false, $false, $true, 1, $false, "13÷7=1, 6", 2) | Out-Null

$c = $t.Cell(17, 3)
$rng = $c.Range
$rng.MoveEnd(1, -1)
$rng.Find.Execute("43÷8=5, 3", $true, $false, $false, $false, $false, $true, 1, $false, "94÷8=11, 6", 2) | Out-Null

$c = $t.Cell(17, 4)
$rng = $c.Range
$rng.MoveEnd(1, -1)
$rng.Find.Execute("35÷8=4, 3", $true, $false, $false, $false, $false, $true, 1, $false, "14÷2=7, 0", 2) | Out-Null

$c = $t.Cell(17, 5)
$rng = $c.Range
$rng.MoveEnd(1, -1)
$rng.Find.Execute("35÷5=7, 0", $true, $false, $false, $false, $false, $true, 1, $false, "62÷8=7, 6", 2) | Out-Null

